$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.748.27"
$ws.Range("E2").Value = "  +8.29%  "
$ws.Range("D3").Value = "1.948.66"
$ws.Range("E3").Value = "  +6.77%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").Value = "'342.06"
$ws.Range("E5").Value = "  +2.97%  "
$ws.Range("D6").Value = "'0.9999"
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("D7").Value = "'0.4779"
$ws.Range("E7").Value = "  +4.19%  "
$ws.Range("D8").Value = "'0.4155"
$ws.Range("E8").Value = "  +8.82%  "
$ws.Range("D9").Value = "'48.32"
$ws.Range("E9").Value = "  +5.50%  "
$ws.Range("E10").Value = "  +5.17%  "
$ws.Range("D11").Value = "'1.044"
$ws.Range("E11").Value = "  +8.74%  "
$ws.Range("D12").Value = "'22.70"
$ws.Range("E12").Value = "  +7.74%  "
$ws.Range("D13").Value = "1.938.28"
$ws.Range("E13").Value = "  +5.39%  "
$ws.Range("D14").Value = "'6.197"
$ws.Range("E14").Value = "  +5.95%  "
$ws.Range("D15").Value = "'7.426"
$ws.Range("E15").Value = "  +4.75%  "
$ws.Range("D16").Value = "'92.44"
$ws.Range("E16").Value = "  +3.26%  "
$ws.Range("D17").Value = "'1.001"
$ws.Range("E17").Value = "  -0.22%  "
$ws.Range("D18").Value = "'0.00001065"
$ws.Range("E18").Value = "  +4.30%  "
$ws.Range("D19").Value = "'0.06670"
$ws.Range("E19").Value = "  +1.22%  "
$ws.Range("D20").Value = "'18.07"
$ws.Range("E20").Value = "  +5.36%  "
$ws.Range("D21").Value = "'0.9999"
$ws.Range("E21").Value = "  -0.27%  "
$ws.Range("D22").Value = "29.712.01"
$ws.Range("E22").Value = "  +8.23%  "
$ws.Range("D23").Value = "'5.612"
$ws.Range("E23").Value = "  +6.00%  "
$ws.Range("D24").Value = "'11.27"
$ws.Range("E24").Value = "  +3.81%  "
$ws.Range("D25").Value = "'2.282"
$ws.Range("E25").Value = "  +1.02%  "
$ws.Range("D26").Value = "2.174.08"
$ws.Range("E26").Value = "  +5.83%  "
$ws.Range("D27").Value = "'160.72"
$ws.Range("E27").Value = "  +2.52%  "
$ws.Range("D28").Value = "'20.22"
$ws.Range("E28").Value = "  +4.54%  "
$ws.Range("D29").Value = "'2.199"
$ws.Range("E29").Value = "  +7.38%  "
$ws.Range("D30").Value = "'5.676"
$ws.Range("E30").Value = "  +7.59%  "
$ws.Range("D31").Value = "'122.48"
$ws.Range("E31").Value = "  +3.85%  "
$ws.Range("D32").Value = "'1.029"
$ws.Range("E32").Value = "  +10.19%  "
$ws.Range("D33").Value = "'0.09638"
$ws.Range("E33").Value = "  +3.52%  "
$ws.Range("D34").Value = "'1.481"
$ws.Range("E34").Value = "  +12.51%  "
$ws.Range("D35").Value = "'3.678"
$ws.Range("E35").Value = "  +3.17%  "
$ws.Range("D36").Value = "'5.490"
$ws.Range("E36").Value = "  +5.02%  "
$ws.Range("D37").Value = "'0.06317"
$ws.Range("E37").Value = "  +6.58%  "
$ws.Range("D38").Value = "'0.02336"
$ws.Range("E38").Value = "  +6.70%  "
$ws.Range("D39").Value = "'8.617"
$ws.Range("E39").Value = "  +6.02%  "
$ws.Range("D40").Value = "'1.198"
$ws.Range("E40").Value = "  +4.79%  "
$ws.Range("D41").Value = "'0.6122"
$ws.Range("E41").Value = "  +6.26%  "
$ws.Range("D42").Value = "'10.74"
$ws.Range("E42").Value = "  +7.84%  "
$ws.Range("D43").Value = "'0.1903"
$ws.Range("E43").Value = "  +4.69%  "
$ws.Range("D44").Value = "'0.9998"
$ws.Range("E44").Value = "  -0.25%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "'2.416"
$ws.Range("E45").Value = "  +34.04%  "
$ws.Range("B46").Value = "WEMIXTOKEN"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "'1.271"
$ws.Range("E46").Value = "  -0.31%  "
$ws.Range("D47").Value = "'12.63"
$ws.Range("E47").Value = "  +6.69%  "
$ws.Range("D48").Value = "'0.5718"
$ws.Range("E48").Value = "  +5.78%  "
$ws.Range("D49").Value = "'2.005"
$ws.Range("E49").Value = "  +6.94%  "
$ws.Range("D50").Value = "'0.07411"
$ws.Range("E50").Value = "  +12.64%  "
$ws.Range("D51").Value = "'114.08"
$ws.Range("E51").Value = "  +3.39%  "
